$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Octubre de 2020 a las 02:45"

# --- Update data for countries whose totals changed but kept their rank/row ---

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 7494178
$ws.Cells.Item(4, 3).Value = 46896
$ws.Cells.Item(4, 4).Value = 4731303
$ws.Cells.Item(4, 5).Value = 2550220
$ws.Cells.Item(4, 7).Value = 915
$ws.Cells.Item(4, 8).Value = 212655

# Row 11: Argentina
$ws.Cells.Item(11, 2).Value = 765002
$ws.Cells.Item(11, 3).Value = 14001
$ws.Cells.Item(11, 4).Value = 603140
$ws.Cells.Item(11, 5).Value = 141574
$ws.Cells.Item(11, 7).Value = 3352
$ws.Cells.Item(11, 8).Value = 20288

# Row 129: Surinam
$ws.Cells.Item(129, 2).Value = 4891
$ws.Cells.Item(129, 3).Value = 14
$ws.Cells.Item(129, 4).Value = 4702
$ws.Cells.Item(129, 5).Value = 84
$ws.Cells.Item(129, 7).Value = 1
$ws.Cells.Item(129, 8).Value = 105

# Row 200: San Bartolome
$ws.Cells.Item(200, 2).Value = 54
$ws.Cells.Item(200, 3).Value = 6
$ws.Cells.Item(200, 4).Value = 37
$ws.Cells.Item(200, 5).Value = 17

# --- Rows that swapped rank/position due to the data refresh ---

# Rows 181/182: San Martin (Parte Francesa) overtakes Curazao and gets new data;
# Curazao keeps its previous (unchanged) data and drops one place.
$ws.Cells.Item(181, 1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(181, 2).Value = 403
$ws.Cells.Item(181, 3).Value = 20
$ws.Cells.Item(181, 4).Value = 309
$ws.Cells.Item(181, 5).Value = 86
$ws.Cells.Item(181, 8).Value = 8

$ws.Cells.Item(182, 1).Value = "Curazao"
$ws.Cells.Item(182, 2).Value = 399
$ws.Cells.Item(182, 3).Value = 7
$ws.Cells.Item(182, 4).Value = 185
$ws.Cells.Item(182, 5).Value = 213
$ws.Cells.Item(182, 8).Value = 1

# Rows 207/208: Santa Lucia and Nueva Caledonia swap places (tied totals)
$ws.Cells.Item(207, 1).Value = "Santa Lucia"
$ws.Cells.Item(208, 1).Value = "Nueva Caledonia"

# Rows 215/216: Islas Malvinas and Montserrat swap places (tied totals),
# each carrying its own previous data along with it.
$ws.Cells.Item(215, 1).Value = "Islas Malvinas"
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 8).Value = 0

$ws.Cells.Item(216, 1).Value = "Montserrat"
$ws.Cells.Item(216, 4).Value = 12
$ws.Cells.Item(216, 8).Value = 1
